# Auto-generated script applying scheduled market-data refresh values
# to the Leve profit tracking tables across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 189.13043
$ws.Range("I33").Value = 193.54546
$ws.Range("K33").Value = 193.54546
$ws.Range("M33").Value = 35.45454000000001

$ws.Range("H70").Value = 1469.2
$ws.Range("I70").Value = 1389
$ws.Range("J70").Value = 1478.1111
$ws.Range("K70").Value = 4167
$ws.Range("L70").Value = 4434.3333
$ws.Range("M70").Value = -3897
$ws.Range("N70").Value = -4974.3333

$ws.Range("H73").Value = 1469.2
$ws.Range("I73").Value = 1389
$ws.Range("J73").Value = 1478.1111
$ws.Range("K73").Value = 4167
$ws.Range("L73").Value = 4434.3333
$ws.Range("M73").Value = -3231
$ws.Range("N73").Value = -6306.3333

$ws.Range("H138").Value = 7778.1714
$ws.Range("I138").Value = 5999.5
$ws.Range("J138").Value = 7885.9697
$ws.Range("K138").Value = 17998.5
$ws.Range("L138").Value = 23657.9091
$ws.Range("M138").Value = -12858.5
$ws.Range("N138").Value = -33937.9091

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3967.85
$ws.Range("I32").Value = 3967.85
$ws.Range("K32").Value = 3967.85
$ws.Range("M32").Value = -3680.85

$ws.Range("H110").Value = 3895.4443
$ws.Range("I110").Value = 3997.1667
$ws.Range("J110").Value = 3692
$ws.Range("K110").Value = 3997.1667
$ws.Range("L110").Value = 3692
$ws.Range("M110").Value = -1952.1667
$ws.Range("N110").Value = -7782

$ws.Range("H122").Value = 2416.9
$ws.Range("I122").Value = 2416.9
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7250.700000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4800.700000000001
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3704.5715
$ws.Range("I132").Value = 2811.6667
$ws.Range("K132").Value = 8435.000100000001
$ws.Range("M132").Value = -5905.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 120.4
$ws.Range("I7").Value = 90
$ws.Range("J7").Value = 166
$ws.Range("K7").Value = 90
$ws.Range("L7").Value = 166
$ws.Range("M7").Value = 23
$ws.Range("N7").Value = -392

$ws.Range("H31").Value = 4924.4814
$ws.Range("I31").Value = 1723.75
$ws.Range("J31").Value = 9580.091
$ws.Range("K31").Value = 1723.75
$ws.Range("L31").Value = 9580.091
$ws.Range("M31").Value = -1428.75
$ws.Range("N31").Value = -10170.091

$ws.Range("H34").Value = 4924.4814
$ws.Range("I34").Value = 1723.75
$ws.Range("J34").Value = 9580.091
$ws.Range("K34").Value = 1723.75
$ws.Range("L34").Value = 9580.091
$ws.Range("M34").Value = -1521.75
$ws.Range("N34").Value = -9984.091

$ws.Range("H97").Value = 30045.5
$ws.Range("J97").Value = 30045.5
$ws.Range("L97").Value = 30045.5
$ws.Range("N97").Value = -32027.5

$ws.Range("H107").Value = 796.3333
$ws.Range("I107").Value = 414.75
$ws.Range("K107").Value = 414.75
$ws.Range("M107").Value = 1505.25

$ws.Range("H130").Value = 66666
$ws.Range("J130").Value = 66666
$ws.Range("L130").Value = 66666
$ws.Range("N130").Value = -76706

$ws.Range("H132").Value = 2735.7368
$ws.Range("I132").Value = 1399.5
$ws.Range("K132").Value = 4198.5
$ws.Range("M132").Value = -1668.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 198980
$ws.Range("J37").Value = 198980
$ws.Range("L37").Value = 596940
$ws.Range("N37").Value = -597164

$ws.Range("H68").Value = 1120.7778
$ws.Range("J68").Value = 1110.875
$ws.Range("L68").Value = 3332.625
$ws.Range("N68").Value = -4954.625

$ws.Range("H71").Value = 1120.7778
$ws.Range("J71").Value = 1110.875
$ws.Range("L71").Value = 9997.875
$ws.Range("N71").Value = -18109.875

$ws.Range("H113").Value = 668.0769
$ws.Range("J113").Value = 733.75
$ws.Range("L113").Value = 2201.25
$ws.Range("N113").Value = -6541.25

$ws.Range("H117").Value = 1912
$ws.Range("I117").Value = 2216.6667
$ws.Range("J117").Value = 998
$ws.Range("K117").Value = 6650.000100000001
$ws.Range("L117").Value = 2994
$ws.Range("M117").Value = -3208.000100000001
$ws.Range("N117").Value = -9878

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5746.1875
$ws.Range("I80").Value = 5367.125
$ws.Range("K80").Value = 5367.125
$ws.Range("M80").Value = -4369.125

$ws.Range("H83").Value = 5746.1875
$ws.Range("I83").Value = 5367.125
$ws.Range("K83").Value = 26835.625
$ws.Range("M83").Value = -21843.625

$ws.Range("H107").Value = 378.33334
$ws.Range("I107").Value = 299
$ws.Range("J107").Value = 418
$ws.Range("K107").Value = 299
$ws.Range("L107").Value = 418
$ws.Range("M107").Value = 1621
$ws.Range("N107").Value = -4258

$ws.Range("H126").Value = 2331.6667
$ws.Range("I126").Value = 2331.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6995.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4525.000100000001
$ws.Range("N126").ClearContents()

$ws.Range("H135").Value = 125000
$ws.Range("J135").Value = 125000
$ws.Range("L135").Value = 125000
$ws.Range("N135").Value = -135140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H7").Value = 3771
$ws.Range("I7").Value = 3157
$ws.Range("J7").Value = 4078
$ws.Range("K7").Value = 3157
$ws.Range("L7").Value = 4078
$ws.Range("M7").Value = -3045
$ws.Range("N7").Value = -4302

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H122").Value = 3994
$ws.Range("I122").Value = 3994
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11982
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9532
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 3771
$ws.Range("I126").Value = 3157
$ws.Range("J126").Value = 4078
$ws.Range("K126").Value = 9471
$ws.Range("L126").Value = 12234
$ws.Range("M126").Value = -7001
$ws.Range("N126").Value = -17174

$ws.Range("H132").Value = 3440.2
$ws.Range("I132").Value = 2161.4
$ws.Range("K132").Value = 6484.200000000001
$ws.Range("M132").Value = -3954.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6499.6665
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 6499.6665
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6499.6665
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7747.6665

$ws.Range("H65").Value = 6499.6665
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 6499.6665
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 32498.3325
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -38738.3325

$ws.Range("H96").Value = 3085
$ws.Range("J96").Value = 1889
$ws.Range("L96").Value = 1889
$ws.Range("N96").Value = -4635

$ws.Range("H113").Value = 1272.25
$ws.Range("I113").Value = 1212.25
$ws.Range("K113").Value = 3636.75
$ws.Range("M113").Value = -1466.75

$ws.Range("H132").Value = 2450.2354
$ws.Range("I132").Value = 2240.96
$ws.Range("K132").Value = 6722.88
$ws.Range("M132").Value = -4192.88

$ws.Range("H136").Value = 5338.433
$ws.Range("I136").Value = 5659.56
$ws.Range("J136").Value = 3732.8
$ws.Range("K136").Value = 16978.68
$ws.Range("L136").Value = 11198.4
$ws.Range("M136").Value = -14428.68

